# Auto-generated edit script: update hashcode values in column B
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hashcode.csv")

$updates = @{
    84 = "46c8b004fe3a12e28ad00319fdb83d45"
    86 = "d9afb6bfd62af8b9d4dae18c3627ed50"
    88 = "6670521f629b500de1ad7a44f8af0ac3"
    134 = "b1d389d724c204a8cdd69e4ef2a4c80a"
    148 = "fcd77afe467c8ebc17f91c5b7733d8a2"
    186 = "cfc31542d110f939cc49ba628e8dade7"
    195 = "199f5bee8b64eef4d3ce639d5ee83a45"
    196 = "f98db8c6bbfa9e154ab69a88e8c28c40"
    214 = "a41ed8fea075d599e76133e0a7f103f2"
    273 = "8c6e2b75376b8490b816902250befb49"
    318 = "86f173682dfc13b946926536be2dd472"
    338 = "9dd383fa4a1dc596a7edfd89dd3f91d6"
    340 = "772d0fdaf33355458672d788f59ddc85"
    348 = "8200bca2df3947a617bcbf7917764f1f"
    357 = "c6f999690ed77fd7c3886fe8d01a8eef"
    359 = "566999730ad2fb72c76567cc59af78ec"
    363 = "bc161f3027a3f00c756aeef058b753a4"
    380 = "c4fd6a80e4a21d17aefa86ca3e38df33"
    390 = "71e4b8478f452dc0aacb692f3ee15615"
    392 = "5dcf754b25f865b537e9ee08c27e372a"
    395 = "0aa42c3e5842964b143da8180d3bb392"
    412 = "cc892fa61f4d68caaf5b8696e8c2248c"
    433 = "d8dcf9f300a60ae547d1469c8cc61e12"
    441 = "db0b37a75d5af39c44c29a5166e7bc1d"
    463 = "ae39cf2b3c81e8771d4d83770cbfdceb"
    487 = "1fe0235aed6099e8f9cae19b2c91ae37"
    491 = "fce320f030c3dd011ba7cada5ffe4f33"
    512 = "482c94f933b0e3dfbd916d15bf78c9a6"
    531 = "31f69821165e8eae7b2dddbefdb0ac95"
    545 = "72fa49926b7c0a331f7c048fa1ce7b9b"
    563 = "ba9d1f54292d1a172ca42c55c7193514"
    607 = "ca08fc086a5116ce02a986c677815ea5"
    660 = "e3846db6488435f56d01f5f0d8029ba0"
    661 = "bc0ba5600326d366610d8e8c4660bc7a"
    691 = "a8aa74d27a1928c88d25058aeafab667"
    717 = "5d91c5df5f813b2d2e03cae38c96e882"
    725 = "05559272def7c987ce6d3c88d83b9f20"
    727 = "909250fe5e6250751566411a7195725c"
    740 = "1e7c73cb11ee72bd7754d2462735e3ec"
    753 = "3ab88564f41e122d4c712b8cd2d17006"
    755 = "16b202040160feaff032553e4e0f7211"
    772 = "cf19049c07b809407be05a9d49a1742f"
    809 = "ad2752e3512e4c597db2460e9dcc930a"
    828 = "47ad0ae47b1bf68b36a8346401da63de"
    831 = "df207d48540e4cac1fdb9eecb685d11a"
    908 = "6a96b3d85850c0e1f8838a8aabb8c092"
    921 = "edd137bbe0b4627809b6bb96afa5e6df"
    930 = "b4f7367829d04896214390c704d14f06"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}

